# Reverse the order of the comma-separated "Recorded By" entries in column G
# for every data row on the "Session Analysis Results" sheet.
#
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#      "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"
#      "System, backup@backdoor.com, system" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }

    [array]::Reverse($parts)
    $newVal = [string]::Join(", ", $parts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
